$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 153; shifts existing rows 153:233 down to 154:234
$ws.Rows("153:153").Insert()

# Populate the newly inserted row 153 with the new data record
$ws.Cells.Item(153, 1).Value = 3
$ws.Cells.Item(153, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(153, 3).Value = "Coquimbo"
$ws.Cells.Item(153, 4).Value = 44460
$ws.Cells.Item(153, 5).Value = 5
$ws.Cells.Item(153, 6).Value = 100112017
$ws.Cells.Item(153, 7).Value = "Apio"
$ws.Cells.Item(153, 8).Value = "Americana (o)"
$ws.Cells.Item(153, 9).Value = "Primera"
$ws.Cells.Item(153, 10).Value = 230
$ws.Cells.Item(153, 11).Value = 9000
$ws.Cells.Item(153, 12).Value = 9500
$ws.Cells.Item(153, 13).Value = 9239
$ws.Cells.Item(153, 14).Value = "`$/docena de matas"
$ws.Cells.Item(153, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(153, 16).Value = 1540
$ws.Cells.Item(153, 17).Value = 6
$ws.Cells.Item(153, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Cells.Item(153, 4).NumberFormat = $ws.Cells.Item(154, 4).NumberFormat
